$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''30.175.88'
$ws.Range("E2").Value = '  +0.30%  '

# Row 3
$ws.Range("D3").Value = '''1.908.66'
$ws.Range("E3").Value = '  -0.26%  '

# Row 4
$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '  +0.21%  '

# Row 5
$ws.Range("D5").Value = '''0.8224'
$ws.Range("E5").Value = '  +4.46%  '

# Row 6
$ws.Range("D6").Value = '''243.34'
$ws.Range("E6").Value = '  -0.07%  '

# Row 7
$ws.Range("E7").Value = '  -0.04%  '

# Row 8
$ws.Range("D8").Value = '''0.3262'
$ws.Range("E8").Value = '  +2.97%  '

# Row 9
$ws.Range("E9").Value = '  +2.24%  '

# Row 10
$ws.Range("D10").Value = '''0.07051'
$ws.Range("E10").Value = '  +1.75%  '

# Row 11
$ws.Range("D11").Value = '''0.08098'
$ws.Range("E11").Value = '  +1.27%  '

# Row 12
$ws.Range("D12").Value = '''0.7647'
$ws.Range("E12").Value = '  +2.36%  '

# Row 13
$ws.Range("D13").Value = '''1.912.40'
$ws.Range("E13").Value = '  -0.01%  '

# Row 14
$ws.Range("D14").Value = '''5.276'
$ws.Range("E14").Value = '  +0.94%  '

# Row 15
$ws.Range("D15").Value = '''92.74'
$ws.Range("E15").Value = '  -0.71%  '

# Row 16
$ws.Range("D16").Value = '''30.176.65'
$ws.Range("E16").Value = '  +0.24%  '

# Row 17
$ws.Range("D17").Value = '''14.19'
$ws.Range("E17").Value = '  +1.36%  '

# Row 18
$ws.Range("D18").Value = '''5.900'
$ws.Range("E18").Value = '  -0.65%  '

# Row 19
$ws.Range("D19").Value = '''245.36'

# Row 20
$ws.Range("D20").Value = '''0.000007784'
$ws.Range("E20").Value = '  +0.03%  '

# Row 21
$ws.Range("D21").Value = '''2.166.04'
$ws.Range("E21").Value = '  +0.17%  '

# Row 22
$ws.Range("E22").Value = '  +0.04%  '

# Row 23
$ws.Range("D23").Value = '''1.002'
$ws.Range("E23").Value = '  +0.04%  '

# Row 24
$ws.Range("D24").Value = '''7.039'
$ws.Range("E24").Value = '  +1.81%  '

# Row 25
$ws.Range("D25").Value = '''0.1677'
$ws.Range("E25").Value = '  +21.56%  '

# Row 26
$ws.Range("D26").Value = '''9.311'
$ws.Range("E26").Value = '  -0.07%  '

# Row 27
$ws.Range("D27").Value = '''166.54'
$ws.Range("E27").Value = '  -1.60%  '

# Row 28
$ws.Range("D28").Value = '''19.00'
$ws.Range("E28").Value = '  +0.35%  '

# Row 29
$ws.Range("D29").Value = '''2.101'
$ws.Range("E29").Value = '  +2.58%  '

# Row 30
$ws.Range("D30").Value = '''1.370'
$ws.Range("E30").Value = '  -0.49%  '

# Row 31
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '''1.525'
$ws.Range("E31").Value = '  -0.10%  '

# Row 32
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '''0.06007'
$ws.Range("E32").Value = '  +4.61%  '

# Row 33
$ws.Range("D33").Value = '''4.301'
$ws.Range("E33").Value = '  -1.02%  '

# Row 34
$ws.Range("D34").Value = '''4.087'
$ws.Range("E34").Value = '  -0.72%  '

# Row 35
$ws.Range("D35").Value = '''1.271'
$ws.Range("E35").Value = '  +0.77%  '

# Row 36
$ws.Range("D36").Value = '''0.7344'
$ws.Range("E36").Value = '  -0.39%  '

# Row 37
$ws.Range("E37").Value = '  -0.50%  '

# Row 38
$ws.Range("D38").Value = '''0.01932'
$ws.Range("E38").Value = '  +0.52%  '

# Row 40
$ws.Range("D40").Value = '''0.4457'
$ws.Range("E40").Value = '  +0.13%  '

# Row 41
$ws.Range("D41").Value = '''73.10'
$ws.Range("E41").Value = '  +0.58%  '

# Row 42
$ws.Range("D42").Value = '''5.962'
$ws.Range("E42").Value = '  -3.39%  '

# Row 43
$ws.Range("D43").Value = '''0.8531'
$ws.Range("E43").Value = '  +1.92%  '

# Row 44
$ws.Range("D44").Value = '''1.002'
$ws.Range("E44").Value = '  -0.11%  '

# Row 45
$ws.Range("D45").Value = '''1.905'
$ws.Range("E45").Value = '  +0.34%  '

# Row 46
$ws.Range("D46").Value = '''102.29'
$ws.Range("E46").Value = '  +1.50%  '

# Row 47
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = '''7.586'
$ws.Range("E47").Value = '  +0.06%  '

# Row 48
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''9.844'
$ws.Range("E48").Value = '  +0.11%  '

# Row 49
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '''1.006.09'
$ws.Range("E49").Value = '  +1.42%  '

# Row 50
$ws.Range("D50").Value = '''2.062.87'
$ws.Range("E50").Value = '  +0.05%  '

# Row 51
$ws.Range("D51").Value = '''1.554'
$ws.Range("E51").Value = '  +3.77%  '
